$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Methaanslip tijdens CO2 vervloeiing uit bedrijf (garantie 2 g CH4/Nm3 biogas)"
$ws.Range("B1").Value = "##"
$ws.Range("C1").Value = "g CH4/Nm3 biogas"
$ws.Range("E1").Value = "Methaanslip tijdens CO2 vervloeiing actief (garantie < 0,01 %)"
$ws.Range("F1").Value = "##"
$ws.Range("G1").Value = "%"
